$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 177, shifting the existing rows 177-254 down to 180-257.
$ws.Range("A177:A179").EntireRow.Insert()

# Row 177: new "Primera" quality record dated 2022-10-07 (serial 44841)
$ws.Cells.Item(177, 1).Value = 12
$ws.Cells.Item(177, 2).Value = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(177, 3).Value = "Metropolitana"
$ws.Cells.Item(177, 4).Value = 44841
$ws.Cells.Item(177, 5).Value = 13
$ws.Cells.Item(177, 6).Value = 100112043
$ws.Cells.Item(177, 7).Value = "Pepino dulce"
$ws.Cells.Item(177, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(177, 9).Value = "Primera"
$ws.Cells.Item(177, 10).Value = 80
$ws.Cells.Item(177, 11).Value = 18000
$ws.Cells.Item(177, 12).Value = 18000
$ws.Cells.Item(177, 13).Value = 18000
$ws.Cells.Item(177, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(177, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(177, 16).Value = 1000
$ws.Cells.Item(177, 17).Value = 18
$ws.Cells.Item(177, 18).Value = "Hortaliza"

# Row 178: new "Segunda" quality record dated 2022-10-07 (serial 44841)
$ws.Cells.Item(178, 1).Value = 12
$ws.Cells.Item(178, 2).Value = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(178, 3).Value = "Metropolitana"
$ws.Cells.Item(178, 4).Value = 44841
$ws.Cells.Item(178, 5).Value = 13
$ws.Cells.Item(178, 6).Value = 100112043
$ws.Cells.Item(178, 7).Value = "Pepino dulce"
$ws.Cells.Item(178, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(178, 9).Value = "Segunda"
$ws.Cells.Item(178, 10).Value = 95
$ws.Cells.Item(178, 11).Value = 14000
$ws.Cells.Item(178, 12).Value = 14000
$ws.Cells.Item(178, 13).Value = 14000
$ws.Cells.Item(178, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(178, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(178, 16).Value = 778
$ws.Cells.Item(178, 17).Value = 18
$ws.Cells.Item(178, 18).Value = "Hortaliza"

# Row 179: new "Tercera" quality record dated 2022-10-07 (serial 44841)
$ws.Cells.Item(179, 1).Value = 12
$ws.Cells.Item(179, 2).Value = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(179, 3).Value = "Metropolitana"
$ws.Cells.Item(179, 4).Value = 44841
$ws.Cells.Item(179, 5).Value = 13
$ws.Cells.Item(179, 6).Value = 100112043
$ws.Cells.Item(179, 7).Value = "Pepino dulce"
$ws.Cells.Item(179, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(179, 9).Value = "Tercera"
$ws.Cells.Item(179, 10).Value = 120
$ws.Cells.Item(179, 11).Value = 10000
$ws.Cells.Item(179, 12).Value = 10000
$ws.Cells.Item(179, 13).Value = 10000
$ws.Cells.Item(179, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(179, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(179, 16).Value = 556
$ws.Cells.Item(179, 17).Value = 18
$ws.Cells.Item(179, 18).Value = "Hortaliza"
